$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema3c"
$ws.Range("C2").Value = "Plxnd1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2583983333333333
$ws.Range("H2").Value = 0.775195
$ws.Range("I2").Value = 0.007195239230717037
$ws.Range("J2").Value = 0.007195239230717038
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 81.979392
$ws.Range("N2").Value = 245.938176
$ws.Range("O2").Value = 0.523851352180617
$ws.Range("P2").Value = 0.523851352180617
$ws.Range("Q2").Value = 21.18333826048
$ws.Range("R2").Value = 190.65004434432
$ws.Range("S2").Value = 0.003769235800274142
$ws.Range("T2").Value = 0.003769235800274143

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema3c"
$ws.Range("C3").Value = "Plxnd1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2583983333333333
$ws.Range("H3").Value = 0.775195
$ws.Range("I3").Value = 0.007195239230717037
$ws.Range("J3").Value = 0.007195239230717038
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 22.17197066666667
$ws.Range("N3").Value = 66.515912
$ws.Range("O3").Value = 0.1416797140218155
$ws.Range("P3").Value = 0.1416797140218155
$ws.Range("Q3").Value = 5.729200266982223
$ws.Range("R3").Value = 51.56280240284
$ws.Range("S3").Value = 0.001019419436526537
$ws.Range("T3").Value = 0.001019419436526538

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema3c"
$ws.Range("C4").Value = "Plxnd1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2583983333333333
$ws.Range("H4").Value = 0.775195
$ws.Range("I4").Value = 0.007195239230717037
$ws.Range("J4").Value = 0.007195239230717038
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 29.98794266666667
$ws.Range("N4").Value = 89.963828
$ws.Range("O4").Value = 0.1916240646801595
$ws.Range("P4").Value = 0.1916240646801595
$ws.Range("Q4").Value = 7.748834405162223
$ws.Range("R4").Value = 69.73950964646001
$ws.Range("S4").Value = 0.001378780987736143
$ws.Range("T4").Value = 0.001378780987736143

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sema3c"
$ws.Range("C5").Value = "Plxnd1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2583983333333333
$ws.Range("H5").Value = 0.775195
$ws.Range("I5").Value = 0.007195239230717037
$ws.Range("J5").Value = 0.007195239230717038
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 22.35431
$ws.Range("N5").Value = 67.06293000000001
$ws.Range("O5").Value = 0.142844869117408
$ws.Range("P5").Value = 0.1428448691174081
$ws.Range("Q5").Value = 5.776316446816668
$ws.Range("R5").Value = 51.98684802135001
$ws.Range("S5").Value = 0.001027803006180215
$ws.Range("T5").Value = 0.001027803006180215

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema3c"
$ws.Range("C6").Value = "Plxnd1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 35.10542566666667
$ws.Range("H6").Value = 105.316277
$ws.Range("I6").Value = 0.9775292770250872
$ws.Range("J6").Value = 0.9775292770250873
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 81.979392
$ws.Range("N6").Value = 245.938176
$ws.Range("O6").Value = 0.523851352180617
$ws.Range("P6").Value = 0.523851352180617
$ws.Range("Q6").Value = 2877.921452054528
$ws.Range("R6").Value = 25901.29306849075
$ws.Range("S6").Value = 0.5120800335657328
$ws.Range("T6").Value = 0.5120800335657328

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema3c"
$ws.Range("C7").Value = "Plxnd1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 35.10542566666667
$ws.Range("H7").Value = 105.316277
$ws.Range("I7").Value = 0.9775292770250872
$ws.Range("J7").Value = 0.9775292770250873
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.17197066666667
$ws.Range("N7").Value = 66.515912
$ws.Range("O7").Value = 0.1416797140218155
$ws.Range("P7").Value = 0.1416797140218155
$ws.Range("Q7").Value = 778.3564681221804
$ws.Range("R7").Value = 7005.208213099624
$ws.Range("S7").Value = 0.1384960684168664
$ws.Range("T7").Value = 0.1384960684168664

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sema3c"
$ws.Range("C8").Value = "Plxnd1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 35.10542566666667
$ws.Range("H8").Value = 105.316277
$ws.Range("I8").Value = 0.9775292770250872
$ws.Range("J8").Value = 0.9775292770250873
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 29.98794266666667
$ws.Range("N8").Value = 89.963828
$ws.Range("O8").Value = 0.1916240646801595
$ws.Range("P8").Value = 0.1916240646801595
$ws.Range("Q8").Value = 1052.739492180929
$ws.Range("R8").Value = 9474.655429628358
$ws.Range("S8").Value = 0.1873181334074049
$ws.Range("T8").Value = 0.1873181334074049

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sema3c"
$ws.Range("C9").Value = "Plxnd1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 35.10542566666667
$ws.Range("H9").Value = 105.316277
$ws.Range("I9").Value = 0.9775292770250872
$ws.Range("J9").Value = 0.9775292770250873
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 22.35431
$ws.Range("N9").Value = 67.06293000000001
$ws.Range("O9").Value = 0.142844869117408
$ws.Range("P9").Value = 0.1428448691174081
$ws.Range("Q9").Value = 784.7575680346234
$ws.Range("R9").Value = 7062.818112311611
$ws.Range("S9").Value = 0.1396350416350831
$ws.Range("T9").Value = 0.1396350416350831

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema3c"
$ws.Range("C10").Value = "Plxnd1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5485793333333333
$ws.Range("H10").Value = 1.645738
$ws.Range("I10").Value = 0.01527548374419571
$ws.Range("J10").Value = 0.01527548374419571
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 81.979392
$ws.Range("N10").Value = 245.938176
$ws.Range("O10").Value = 0.523851352180617
$ws.Range("P10").Value = 0.523851352180617
$ws.Range("Q10").Value = 44.972200210432
$ws.Range("R10").Value = 404.749801893888
$ws.Range("S10").Value = 0.008002082814609958
$ws.Range("T10").Value = 0.008002082814609958

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Sema3c"
$ws.Range("C11").Value = "Plxnd1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5485793333333333
$ws.Range("H11").Value = 1.645738
$ws.Range("I11").Value = 0.01527548374419571
$ws.Range("J11").Value = 0.01527548374419571
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 22.17197066666667
$ws.Range("N11").Value = 66.515912
$ws.Range("O11").Value = 0.1416797140218155
$ws.Range("P11").Value = 0.1416797140218155
$ws.Range("Q11").Value = 12.16308488700622
$ws.Range("R11").Value = 109.467763983056
$ws.Range("S11").Value = 0.00216422616842254
$ws.Range("T11").Value = 0.002164226168422541

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Sema3c"
$ws.Range("C12").Value = "Plxnd1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5485793333333333
$ws.Range("H12").Value = 1.645738
$ws.Range("I12").Value = 0.01527548374419571
$ws.Range("J12").Value = 0.01527548374419571
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 29.98794266666667
$ws.Range("N12").Value = 89.963828
$ws.Range("O12").Value = 0.1916240646801595
$ws.Range("P12").Value = 0.1916240646801595
$ws.Range("Q12").Value = 16.45076559611822
$ws.Range("R12").Value = 148.056890365064
$ws.Range("S12").Value = 0.002927150285018484
$ws.Range("T12").Value = 0.002927150285018485

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Sema3c"
$ws.Range("C13").Value = "Plxnd1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5485793333333333
$ws.Range("H13").Value = 1.645738
$ws.Range("I13").Value = 0.01527548374419571
$ws.Range("J13").Value = 0.01527548374419571
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 22.35431
$ws.Range("N13").Value = 67.06293000000001
$ws.Range("O13").Value = 0.142844869117408
$ws.Range("P13").Value = 0.1428448691174081
$ws.Range("Q13").Value = 12.26311247692667
$ws.Range("R13").Value = 110.36801229234
$ws.Range("S13").Value = 0.002182024476144731
$ws.Range("T13").Value = 0.002182024476144732
